$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Background: the "Errors" localisation sheet has three columns:
#   A = STR_key (english error key), B = INT_code, C = STR_message (chinese)
# Row 174 used to be the single error "playerAlreadyBindGCAccountId" /
# "账号GameCenter账号已经绑定" (code 673).
#
# The account system now distinguishes "already bound" from "not bound":
#   - row 174 keeps code 673 / message "账号GameCenter账号已经绑定", but the
#     key is renamed to "playerAlreadyBindGCAId"
#   - row 175 (new, code 674) is "theGCIdAlreadyHasDatas" /
#     "此GameCenter账号已有玩家数据"
#   - row 176 (new, code 675) is "theGCAccountDoNotHasData" /
#     "此GameCenter账号下无玩家数据"
#
# The write order below is chosen deliberately so the shared-string table
# ends up with the new strings appended in the same order as upstream:
#   ...,347:"账号GameCenter账号已经绑定",
#   348:"此GameCenter账号已有玩家数据", 349:"此GameCenter账号下无玩家数据",
#   350:"playerAlreadyBindGCAId", 351:"theGCIdAlreadyHasDatas",
#   352:"theGCAccountDoNotHasData"
# ---------------------------------------------------------------------------

# 1) Seed the two new Chinese messages while A174 still solely owns the old
#    shared string slot, so the first one reuses it in place; C175 picks up
#    a second reference to the same text so the slot "locks" before we move
#    A174's key off of it.
$ws.Cells.Item(174, 1).Value = "此GameCenter账号已有玩家数据"
$ws.Cells.Item(175, 3).Value = "此GameCenter账号已有玩家数据"
$ws.Cells.Item(176, 3).Value = "此GameCenter账号下无玩家数据"

# 2) Now repoint the keys - these become fresh shared strings.
$ws.Cells.Item(174, 1).Value = "playerAlreadyBindGCAId"
$ws.Cells.Item(175, 1).Value = "theGCIdAlreadyHasDatas"
$ws.Cells.Item(176, 1).Value = "theGCAccountDoNotHasData"

# 3) Error codes for the two new rows.
$ws.Cells.Item(175, 2).Value = 674
$ws.Cells.Item(176, 2).Value = 675

# 4) Match row formatting (20pt custom height, same as every other row).
$ws.Rows.Item(175).RowHeight = 20
$ws.Rows.Item(176).RowHeight = 20

# 5) Leave the selection where the author ended up after adding the rows.
[void]$ws.Range("A177").Select()
